# Update the LR-pairs TPM-derived statistics on the active sheet with
# newly recomputed values (new TPM run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FAPs -> ECs)
$ws.Range("M2").Value = 0.353476
$ws.Range("N2").Value = 1.060428
$ws.Range("O2").Value = 0.07277959798226569
$ws.Range("P2").Value = 0.07277959798226567
$ws.Range("Q2").Value = 3.510247264177333
$ws.Range("R2").Value = 31.592225377596
$ws.Range("S2").Value = 0.07040363036080503
$ws.Range("T2").Value = 0.07040363036080502

# Row 3 (FAPs -> FAPs)
$ws.Range("O3").Value = 0.2796129441040817
$ws.Range("P3").Value = 0.2796129441040817
$ws.Range("S3").Value = 0.2704846812371383
$ws.Range("T3").Value = 0.2704846812371383

# Row 4 (FAPs -> MuSCs)
$ws.Range("O4").Value = 0.6476074579136527
$ws.Range("P4").Value = 0.6476074579136526
$ws.Range("S4").Value = 0.626465621546348
$ws.Range("T4").Value = 0.6264656215463479

# Row 5 (MuSCs -> ECs)
$ws.Range("G5").Value = 0.3351376666666667
$ws.Range("I5").Value = 0.03264606685570879
$ws.Range("M5").Value = 0.353476
$ws.Range("N5").Value = 1.060428
$ws.Range("O5").Value = 0.07277959798226569
$ws.Range("P5").Value = 0.07277959798226567
$ws.Range("Q5").Value = 0.1184631218626667
$ws.Range("R5").Value = 1.066168096764
$ws.Range("S5").Value = 0.002375967621460654
$ws.Range("T5").Value = 0.002375967621460653

# Row 6 (MuSCs -> FAPs)
$ws.Range("G6").Value = 0.3351376666666667
$ws.Range("I6").Value = 0.03264606685570879
$ws.Range("O6").Value = 0.2796129441040817
$ws.Range("P6").Value = 0.2796129441040817
$ws.Range("Q6").Value = 0.455125106349889
$ws.Range("R6").Value = 4.096125957149001
$ws.Range("S6").Value = 0.009128262866943416
$ws.Range("T6").Value = 0.009128262866943413

# Row 7 (MuSCs -> MuSCs)
$ws.Range("G7").Value = 0.3351376666666667
$ws.Range("I7").Value = 0.03264606685570879
$ws.Range("O7").Value = 0.6476074579136527
$ws.Range("P7").Value = 0.6476074579136526
$ws.Range("R7").Value = 9.486977532113002
$ws.Range("S7").Value = 0.02114183636730472
$ws.Range("T7").Value = 0.02114183636730472
